$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 for the RF classifier; this shifts the old row 6 ("Ensemble") down to row 7
$ws.Rows.Item(6).Insert()

# Copy formatting from A5 (existing label cell) onto the new A6 label cell so styling matches
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 2: kNN
$ws.Range("A2").Value = "kNN"
$row2 = New-Object 'object[,]' 1,24
$row2[0,0] = 0.650354609929078
$row2[0,1] = 0.7195949511402135
$row2[0,2] = 0.650354609929078
$row2[0,3] = 0.6499308033204485
$row2[0,4] = 0.6802791123312744
$row2[0,5] = 0.7337393381938895
$row2[0,6] = 0.6802791123312744
$row2[0,7] = 0.6814570828139329
$row2[0,8] = 0.6803477465110959
$row2[0,9] = 0.7228110866201387
$row2[0,10] = 0.6803477465110959
$row2[0,11] = 0.6868632838916314
$row2[0,12] = 0.7896591169068863
$row2[0,13] = 0.8032963687462733
$row2[0,14] = 0.7896591169068863
$row2[0,15] = 0.7890163061484834
$row2[0,16] = 0.7704644246167925
$row2[0,17] = 0.7813941049778588
$row2[0,18] = 0.7704644246167925
$row2[0,19] = 0.7701386303190259
$row2[0,20] = 0.7768245252802564
$row2[0,21] = 0.7961989291534056
$row2[0,22] = 0.7768245252802564
$row2[0,23] = 0.7747603453967136
$ws.Range("B2:Y2").Value = $row2

# Row 3: SVM
$ws.Range("A3").Value = "SVM"
$row3 = New-Object 'object[,]' 1,24
$row3[0,0] = 0.8132921528254405
$row3[0,1] = 0.821120947099177
$row3[0,2] = 0.8132921528254405
$row3[0,3] = 0.814105879425558
$row3[0,4] = 0.8154884465797301
$row3[0,5] = 0.825673924973134
$row3[0,6] = 0.8154884465797301
$row3[0,7] = 0.8165697171120989
$row3[0,8] = 0.7876687256920614
$row3[0,9] = 0.8000361554656141
$row3[0,10] = 0.7876687256920614
$row3[0,11] = 0.7882568063073884
$row3[0,12] = 0.8626401281171358
$row3[0,13] = 0.8667069263563292
$row3[0,14] = 0.8626401281171358
$row3[0,15] = 0.86215461177632
$row3[0,16] = 0.8669640814458933
$row3[0,17] = 0.8718403436477388
$row3[0,18] = 0.8669640814458933
$row3[0,19] = 0.8669515214950227
$row3[0,20] = 0.8626630061770761
$row3[0,21] = 0.867521862905375
$row3[0,22] = 0.8626630061770761
$row3[0,23] = 0.8625959207273046
$ws.Range("B3:Y3").Value = $row3

# Row 4: LR
$ws.Range("A4").Value = "LR"
$row4 = New-Object 'object[,]' 1,24
$row4[0,0] = 0.8476549988560971
$row4[0,1] = 0.8543875929254267
$row4[0,2] = 0.8476549988560971
$row4[0,3] = 0.8484254766264593
$row4[0,4] = 0.8712651567147105
$row4[0,5] = 0.875819249782422
$row4[0,6] = 0.8712651567147105
$row4[0,7] = 0.8715984140059622
$row4[0,8] = 0.8541066117593228
$row4[0,9] = 0.8607190002137051
$row4[0,10] = 0.8541066117593228
$row4[0,11] = 0.8548517114911736
$row4[0,12] = 0.8584305650880806
$row4[0,13] = 0.8631187064716158
$row4[0,14] = 0.8584305650880806
$row4[0,15] = 0.858065316010989
$row4[0,16] = 0.8648821779913064
$row4[0,17] = 0.8704746993327854
$row4[0,18] = 0.8648821779913064
$row4[0,19] = 0.8647697933549733
$row4[0,20] = 0.8627316403568978
$row4[0,21] = 0.8680837241827162
$row4[0,22] = 0.8627316403568978
$row4[0,23] = 0.8625653863895216
$ws.Range("B4:Y4").Value = $row4

# Row 5: NB
$ws.Range("A5").Value = "NB"
$row5 = New-Object 'object[,]' 1,24
$row5[0,0] = 0.847609242736216
$row5[0,1] = 0.8555459572210221
$row5[0,2] = 0.847609242736216
$row5[0,3] = 0.8472833278025794
$row5[0,4] = 0.860489590482727
$row5[0,5] = 0.8677771628369717
$row5[0,6] = 0.860489590482727
$row5[0,7] = 0.8604651374710626
$row5[0,8] = 0
$row5[0,9] = 0
$row5[0,10] = 0
$row5[0,11] = 0
$row5[0,12] = 0.828231525966598
$row5[0,13] = 0.8403736679872938
$row5[0,14] = 0.828231525966598
$row5[0,15] = 0.825641043605575
$row5[0,16] = 0.838984214138641
$row5[0,17] = 0.8495676618865732
$row5[0,18] = 0.838984214138641
$row5[0,19] = 0.8367526725646955
$row5[0,20] = 0
$row5[0,21] = 0
$row5[0,22] = 0
$row5[0,23] = 0
$ws.Range("B5:Y5").Value = $row5

# Row 6: RF
$ws.Range("A6").Value = "RF"
$row6 = New-Object 'object[,]' 1,24
$row6[0,0] = 0.8241134751773049
$row6[0,1] = 0.8354763150338858
$row6[0,2] = 0.8241134751773049
$row6[0,3] = 0.8262126705296708
$row6[0,4] = 0.8026309768931595
$row6[0,5] = 0.8107913545846541
$row6[0,6] = 0.8026309768931595
$row6[0,7] = 0.8033751893513579
$row6[0,8] = 0.7661404712880348
$row6[0,9] = 0.7735463560335325
$row6[0,10] = 0.7661404712880348
$row6[0,11] = 0.7648212742271767
$row6[0,12] = 0.8005033173186915
$row6[0,13] = 0.8135117729731988
$row6[0,14] = 0.8005033173186915
$row6[0,15] = 0.8010425954848264
$row6[0,16] = 0.8154655685197895
$row6[0,17] = 0.8268603534401769
$row6[0,18] = 0.8154655685197895
$row6[0,19] = 0.8167799318676142
$row6[0,20] = 0.8155342026996111
$row6[0,21] = 0.8196573393928437
$row6[0,22] = 0.8155342026996111
$row6[0,23] = 0.814054527557275
$ws.Range("B6:Y6").Value = $row6

# Row 7: Ensemble
$ws.Range("A7").Value = "Ensemble"
$row7 = New-Object 'object[,]' 1,24
$row7[0,0] = 0.8583848089681995
$row7[0,1] = 0.8624917545862818
$row7[0,2] = 0.8583848089681995
$row7[0,3] = 0.8586828892117225
$row7[0,4] = 0.8670098375657744
$row7[0,5] = 0.8695623172171265
$row7[0,6] = 0.8670098375657744
$row7[0,7] = 0.8670148878123612
$row7[0,8] = 0.8283687943262411
$row7[0,9] = 0.8402276748290086
$row7[0,10] = 0.8283687943262411
$row7[0,11] = 0.8303815827800278
$row7[0,12] = 0.879775795012583
$row7[0,13] = 0.8832435009377537
$row7[0,14] = 0.879775795012583
$row7[0,15] = 0.8795901069868443
$row7[0,16] = 0.8627316403568978
$row7[0,17] = 0.8680984788891563
$row7[0,18] = 0.8627316403568978
$row7[0,19] = 0.8620746383400251
$row7[0,20] = 0.8625943719972546
$row7[0,21] = 0.8685539188536253
$row7[0,22] = 0.8625943719972546
$row7[0,23] = 0.8628632792574497
$ws.Range("B7:Y7").Value = $row7
